# Daily attendance processing - 2025-11-12 20:48:40
# Normalizes the "Recorded By" column (G): for every session row whose
# recorder list contains more than one entry, rotate the list left by one
# position (the first recorder moves to the end of the list).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "

        if ($parts.Count -gt 1) {
            $rotated = ($parts[1..($parts.Count - 1)] + $parts[0]) -join ", "
            $cell.Value = $rotated
        }
    }
}
